# Update crypto price/volume data per upstream refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.486.09'
$ws.Range('E2').Value = '  -1.11%  '
$ws.Range('D3').Value = '3.105.51'
$ws.Range('E3').Value = '  +1.04%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '556.35'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.74%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '138.40'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -2.36%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '3.098.99'
$ws.Range('E8').Value = '  +1.06%  '
$ws.Range('E9').Value = '  +1.79%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '6.75'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +2.72%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.160'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +4.95%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.456'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +1.23%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '35.38'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -0.93%  '
$ws.Range('E14').Value = '  +0.83%  '
$ws.Range('D15').Value = '3.608.06'
$ws.Range('E15').Value = '  +1.22%  '
$ws.Range('D16').Value = '63.535.82'
$ws.Range('E16').Value = '  -1.09%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.111'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +0.45%  '
$ws.Range('D18').Value = '3.108.50'
$ws.Range('E18').Value = '  +1.02%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '507.64'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +4.25%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '6.72'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +1.73%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '13.71'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.86%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.713'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +4.31%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '7.38'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +2.28%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '12.46'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +0.64%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '78.11'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.47%  '
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('E27').Value = '  +2.79%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '8.26'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -1.00%  '
$ws.Range('E29').Value = '  -1.04%  '
$ws.Range('E30').Value = '  +0.07%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '26.36'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +2.61%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '2.53'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -3.80%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.13'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -0.87%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '540.64'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -8.34%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '58.98'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +13.02%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.93'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +0.25%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '5.20'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -3.05%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.0415'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +4.04%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0803'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +1.96%  '
$ws.Range('D40').Value = '3.095.96'
$ws.Range('E40').Value = '  +3.82%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.118'
$ws.Range('D41').Style = "Normal"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '8.14'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -0.25%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.65'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -6.82%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.258'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +5.73%  '
$ws.Range('B45').Value = 'USDe'
$ws.Range('C45').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.00'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +0.04%  '
$ws.Range('B46').Value = 'Fetch.AI'
$ws.Range('C46').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.11'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +0.97%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '121.60'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +2.34%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '24.14'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -3.78%  '
$ws.Range('E49').Value = '  +0.39%  '
$ws.Range('D50').Value = '0.0₃0502'
$ws.Range('E50').Value = '  -5.33%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.37'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +68.07%  '
